$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vendor name in B2 (shared string "XYZ Ltd" -> "efkjkjfwek  Ltd")
$ws.Range("B2").Value = "efkjkjfwek  Ltd"

# Update invoice number in A2 from text "INV-2024-004" to numeric 111111
$ws.Range("A2").Value = 111111

# Update the active selection to C13
$ws.Range("C13").Select()
